$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.08
$ws.Range("D2").Value = 16.23
$ws.Range("E2").Value = 26.62
$ws.Range("F2").Value = 2.12

# Row 3
$ws.Range("C3").Value = 0.15
$ws.Range("D3").Value = 6.69
$ws.Range("E3").Value = 94.81
$ws.Range("F3").Value = 2.58

# Row 4
$ws.Range("D4").Value = 19.17
$ws.Range("E4").Value = 8.38
$ws.Range("F4").Value = 2.89

# Row 5
$ws.Range("C5").Value = 0.07
$ws.Range("D5").Value = 16.36
$ws.Range("E5").Value = 44.44
$ws.Range("F5").Value = 2.57

# Row 6
$ws.Range("D6").Value = 15.11
$ws.Range("E6").Value = 95
$ws.Range("F6").Value = 1.68

# Row 7
$ws.Range("C7").Value = 0.19
$ws.Range("D7").Value = 16.25
$ws.Range("E7").Value = 146.76
$ws.Range("F7").Value = 17

# Row 8
$ws.Range("C8").Value = 0.25
$ws.Range("D8").Value = 7.58
$ws.Range("E8").Value = 288.59
$ws.Range("F8").Value = 20.98

# Row 9
$ws.Range("C9").Value = 0.07
$ws.Range("E9").Value = 44.84
$ws.Range("F9").Value = 23.08

# Row 10
$ws.Range("C10").Value = 0.16
$ws.Range("D10").Value = 16.38
$ws.Range("E10").Value = 236.49
$ws.Range("F10").Value = 20.27

# Row 11
$ws.Range("C11").Value = 0.22
$ws.Range("D11").Value = 14.94
$ws.Range("E11").Value = 283.33
$ws.Range("F11").Value = 5.58
